# Updated capital structure database
# - Serbia / Bank (Money Center): refreshed the capital-structure figures
#   for the existing two rows (Komercijalna banka's peer row + Komercijalna
#   banka itself), and added a new company row for ALTA banka a.d.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----------------------------------------------------------------------
# Row 2 (existing row, company label stays the placeholder index "2")
# ----------------------------------------------------------------------
# B2 holds a numeric-looking label ("2") that must stay TEXT, like the
# other company-name cells in column B - force text format before typing
# it, then drop the temporary number format again so no stray style sticks.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2"
$ws.Range("B2").Style = "Normal"

$ws.Range("D2").Value = 0.105
$ws.Range("E2").ClearContents()
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 57.047
$ws.Range("L2").Value = 0.3016125621232949
$ws.Range("U2").Value = 36.7
$ws.Range("V2").Value = 0.06801334321719793
$ws.Range("W2").Value = 0.05285284480190647
$ws.Range("X2").Value = 0.04891653365565385
$ws.Range("Y2").Value = 0.003936311146252619
$ws.Range("Z2").Value = 0.2744851046044535
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.04894429125547942
$ws.Range("AC2").Value = -0.04894429125547942
$ws.Range("AD2").Value = 0.255
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0.255
$ws.Range("AG2").Value = -36.445
$ws.Range("AH2").Value = 0.0004723490566911485
$ws.Range("AI2").Value = 0.000315572578599229
$ws.Range("AJ2").Value = -0.07243294809750474
$ws.Range("AK2").Value = -0.047248024580122
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# ----------------------------------------------------------------------
# Row 3 (Komercijalna banka a.d.)
# ----------------------------------------------------------------------
$ws.Range("D3").Value = 0.105
$ws.Range("E3").ClearContents()
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 56.2
$ws.Range("L3").Value = 0.308283049917718
$ws.Range("W3").Value = 0.08370568960381293
$ws.Range("X3").Value = 0.04863235737362938
$ws.Range("Y3").Value = 0.03507333223018355
$ws.Range("Z3").Value = 0.2715221924337206
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0.04863235737362938
$ws.Range("AC3").Value = -0.04863235737362938
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 0
$ws.Range("AK3").Value = 0
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()

# ----------------------------------------------------------------------
# Row 4 (new row: ALTA banka a.d.)
# ----------------------------------------------------------------------
$ws.Range("A4").Value = "Serbia"
$ws.Range("B4").Value = "ALTA banka a.d. (BELEX:JMBN)"
$ws.Range("C4").Value = "Bank (Money Center)"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.847
$ws.Range("L4").Value = 0.1238304093567252
$ws.Range("M4").Value = -0
$ws.Range("N4").Value = -0
$ws.Range("O4").Value = -0
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("S4").Value = 0
$ws.Range("U4").Value = 36.7
$ws.Range("V4").Value = 2.446666666666667
$ws.Range("W4").Value = 0.022
$ws.Range("X4").Value = 0.04920070993767832
$ws.Range("Y4").Value = -0.02720070993767832
$ws.Range("Z4").Value = 0.3870529651425985
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.04925622513732945
$ws.Range("AC4").Value = -0.04925622513732945
$ws.Range("AD4").Value = 0.255
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.255
$ws.Range("AG4").Value = -36.445
$ws.Range("AH4").Value = 0.01671583087512291
$ws.Range("AI4").Value = 0.006495987772258311
$ws.Range("AJ4").Value = 1.699463744462579
$ws.Range("AK4").Value = -14.26418786692759
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
